$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 400, pushing the existing rows 400-422
# down to 402-424. The content that used to live in rows 400-422 ends up
# unchanged (just relocated) in rows 402-424.
$ws.Rows("400:401").Insert()

# The newly inserted rows 400-401 are blank; the old data that used to be
# the most recent observations (rows 398-399, dated 2021-03-15) gets
# pushed further down the weekly history into these new rows, unchanged.
$ws.Range("A398:T399").Copy()
$ws.Range("A400").PasteSpecial()

# Rows 398 and 399 then get the new weekly data (a later observation date
# with different volume/price figures) on top of what was copied down.
$ws.Range("D398").Value = 44516
$ws.Range("M398").Value = 250
$ws.Range("N398").Value = 15000
$ws.Range("O398").Value = 15000
$ws.Range("P398").Value = 15000
$ws.Range("S398").Value = 750

$ws.Range("D399").Value = 44516
$ws.Range("M399").Value = 600
$ws.Range("N399").Value = 14000
$ws.Range("O399").Value = 14000
$ws.Range("P399").Value = 14000
$ws.Range("S399").Value = 700
